$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.511.79'
$ws.Range("E2").Value = '  -0.71%  '
$ws.Range("D3").Value = '3.897.74'
$ws.Range("E3").Value = '  +4.02%  '
$ws.Range("E4").Value = '  +0.22%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '601.66'
$c.Style = "Normal"

$ws.Range("E5").Value = '  -0.11%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '164.52'
$c.Style = "Normal"

$ws.Range("E6").Value = '  -0.74%  '
$ws.Range("D7").Value = '3.900.43'
$ws.Range("E7").Value = '  +4.15%  '
$ws.Range("E8").Value = '  -0.16%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.528'
$c.Style = "Normal"

$ws.Range("E9").Value = '  -1.86%  '
$ws.Range("E10").Value = '  -3.94%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '6.35'
$c.Style = "Normal"

$ws.Range("E11").Value = '  -0.08%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.459'
$c.Style = "Normal"

$ws.Range("E12").Value = '  +0.11%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '36.78'
$c.Style = "Normal"

$ws.Range("E13").Value = '  -2.31%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '0.0000244'
$c.Style = "Normal"

$ws.Range("E14").Value = '  -1.23%  '
$ws.Range("D15").Value = '4.563.29'
$ws.Range("E15").Value = '  +4.35%  '
$ws.Range("D16").Value = '3.941.82'
$ws.Range("E16").Value = '  +4.97%  '
$ws.Range("D17").Value = '68.818.72'
$ws.Range("E17").Value = '  -0.16%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '7.42'
$c.Style = "Normal"

$ws.Range("E18").Value = '  -0.01%  '
$ws.Range("E19").Value = '  -0.88%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '16.92'
$c.Style = "Normal"

$ws.Range("E20").Value = '  -5.40%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '11.13'
$c.Style = "Normal"

$ws.Range("E21").Value = '  -2.59%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '484.38'
$c.Style = "Normal"

$ws.Range("E22").Value = '  -1.08%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '0.716'
$c.Style = "Normal"

$ws.Range("E23").Value = '  -1.05%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '0.0000167'
$c.Style = "Normal"

$ws.Range("E24").Value = '  +12.25%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '84.26'
$c.Style = "Normal"

$ws.Range("E25").Value = '  -0.26%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '2.24'
$c.Style = "Normal"

$ws.Range("E26").Value = '  -1.29%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '12.00'
$c.Style = "Normal"

$ws.Range("E27").Value = '  -2.14%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '10.05'
$c.Style = "Normal"

$ws.Range("E28").Value = '  +0.20%  '
$ws.Range("E29").Value = '  -0.01%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '2.93'
$c.Style = "Normal"

$ws.Range("E30").Value = '  -1.03%  '
$ws.Range("D31").Value = '4.059.99'
$ws.Range("E31").Value = '  +4.35%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '7.81'
$c.Style = "Normal"

$ws.Range("E32").Value = '  -3.93%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '2.37'
$c.Style = "Normal"

$ws.Range("E33").Value = '  -2.26%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '31.80'
$c.Style = "Normal"

$ws.Range("E34").Value = '  +0.76%  '
$ws.Range("D35").Value = '3.851.07'
$ws.Range("E35").Value = '  +4.51%  '
$ws.Range("E36").Value = '  -0.78%  '
$ws.Range("E37").Value = '  +2.77%  '
$ws.Range("E38").Value = '  +0.67%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '5.86'
$c.Style = "Normal"

$ws.Range("E39").Value = '  -0.95%  '
$ws.Range("E40").Value = '  +0.05%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.318'
$c.Style = "Normal"

$ws.Range("E41").Value = '  -1.89%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '435.57'
$c.Style = "Normal"

$ws.Range("E42").Value = '  +2.68%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '2.98'
$c.Style = "Normal"

$ws.Range("E43").Value = '  -3.73%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '48.42'
$c.Style = "Normal"

$ws.Range("E44").Value = '  -0.20%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '1.98'
$c.Style = "Normal"

$ws.Range("E45").Value = '  -0.51%  '
$ws.Range("E46").Value = '  +0.01%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '8.44'
$c.Style = "Normal"

$ws.Range("E47").Value = '  +0.21%  '
$ws.Range("D48").Value = '2.829.73'
$ws.Range("E48").Value = '  +1.68%  '
$ws.Range("B49").Value = 'Monero'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '142.12'
$c.Style = "Normal"

$ws.Range("E49").Value = '  -0.19%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '26.16'
$c.Style = "Normal"

$ws.Range("E50").Value = '  +10.40%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.0354'
$c.Style = "Normal"

$ws.Range("E51").Value = '  +0.54%  '
